# Rule edit - Mead Exclusive Flood Control Space
# Edited execution constraint to check if FC slot was greater than 0 (not just Not NaN).
# This rule was setting Mead outflow during very low scenarios in which flood control
# rules should not be executing.
#
# Net effect on this data workbook: Trace1 (column B) values for rows 2-37 are
# replaced with the corresponding Trace5 (column F) values on every worksheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $srcRange = $ws.Range("F2:F37")
    $dstRange = $ws.Range("B2:B37")
    $dstRange.Value2 = $srcRange.Value2
}
